# Fruta / hortaliza, semanal
# Insert a new weekly record as row 53 in the "Tuna" sheet, pushing the
# existing rows (old 53..116) down by one (new 54..117).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 53..116 down to 54..117, leaving row 53 blank for the new record.
$ws.Rows(53).Insert()

# Populate the new row 53 with this week's data. Columns A,B,C,E,F,G,H,I,J,K
# are identical for every record in this sheet (Tuna / Vega Modelo de
# Temuco); L..T hold the quality/volume/price data for the new week.
$ws.Range("A53").Value = 10
$ws.Range("B53").Value = "Vega Modelo de Temuco"
$ws.Range("C53").Value = "La Araucanía"
$ws.Range("D53").Value = 45225
$ws.Range("E53").Value = 9
$ws.Range("F53").Value = "Fruta"
$ws.Range("G53").Value = 100107
$ws.Range("H53").Value = "Otros"
$ws.Range("I53").Value = 100107011
$ws.Range("J53").Value = "Tuna"
$ws.Range("K53").Value = "Sin especificar"
$ws.Range("L53").Value = "Primera"
$ws.Range("M53").Value = 40
$ws.Range("N53").Value = 38000
$ws.Range("O53").Value = 38000
$ws.Range("P53").Value = 38000
$ws.Range("Q53").Value = "$/caja 16 kilos"
$ws.Range("R53").Value = "Provincia de Los Andes"
$ws.Range("S53").Value = 2375
$ws.Range("T53").Value = 16
